# Updated cryptos list data (Price + Volume(1h) columns, and two swapped rows)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '64.141.26'
$ws.Range("E2").Value = '  -0.22%  '

$ws.Range("D3").Value = "'" + '3.477.73'
$ws.Range("E3").Value = '  -0.69%  '

$ws.Range("E4").Value = '  -0.14%  '

$ws.Range("D5").Value = "'" + '584.80'
$ws.Range("E5").Value = '  -0.27%  '

$ws.Range("D6").Value = "'" + '131.51'
$ws.Range("E6").Value = '  -2.10%  '

$ws.Range("D8").Value = "'" + '0.482'
$ws.Range("E8").Value = '  -0.94%  '

$ws.Range("D9").Value = "'" + '7.62'
$ws.Range("E9").Value = '  +5.04%  '

$ws.Range("D10").Value = "'" + '0.123'
$ws.Range("E10").Value = '  -1.42%  '

$ws.Range("D11").Value = "'" + '0.388'
$ws.Range("E11").Value = '  +0.19%  '

$ws.Range("D12").Value = "'" + '4.063.05'
$ws.Range("E12").Value = '  -0.90%  '

$ws.Range("D13").Value = "'" + '0.120'
$ws.Range("E13").Value = '  -0.15%  '

$ws.Range("E14").Value = '  -2.81%  '

$ws.Range("D15").Value = "'" + '3.470.44'
$ws.Range("E15").Value = '  -1.04%  '

$ws.Range("D16").Value = "'" + '64.101.49'
$ws.Range("E16").Value = '  -0.31%  '

$ws.Range("E17").Value = '  -6.30%  '

$ws.Range("D18").Value = "'" + '9.98'
$ws.Range("E18").Value = '  +0.62%  '

$ws.Range("D19").Value = "'" + '5.69'
$ws.Range("E19").Value = '  -0.98%  '

$ws.Range("D20").Value = "'" + '13.44'
$ws.Range("E20").Value = '  -1.74%  '

$ws.Range("D21").Value = "'" + '385.02'
$ws.Range("E21").Value = '  -2.08%  '

$ws.Range("D22").Value = "'" + '0.570'
$ws.Range("E22").Value = '  -0.23%  '

$ws.Range("D23").Value = "'" + '3.617.18'
$ws.Range("E23").Value = '  -0.72%  '

$ws.Range("D24").Value = "'" + '74.82'
$ws.Range("E24").Value = '  +0.74%  '

$ws.Range("E25").Value = '  +0.11%  '

$ws.Range("E26").Value = '  +0.16%  '

$ws.Range("D27").Value = "'" + '0.0000112'
$ws.Range("E27").Value = '  -2.75%  '

$ws.Range("E28").Value = '  -0.02%  '

$ws.Range("D29").Value = "'" + '2.23'
$ws.Range("E29").Value = '  -0.28%  '

$ws.Range("E30").Value = '  -4.74%  '

$ws.Range("E31").Value = '  -4.27%  '

$ws.Range("D32").Value = "'" + '7.92'
$ws.Range("E32").Value = '  -4.18%  '

$ws.Range("D33").Value = "'" + '3.507.67'
$ws.Range("E33").Value = '  -0.47%  '

$ws.Range("E34").Value = '  +1.17%  '

$ws.Range("D36").Value = "'" + '22.93'
$ws.Range("E36").Value = '  -2.29%  '

$ws.Range("D37").Value = "'" + '5.21'
$ws.Range("E37").Value = '  +0.52%  '

$ws.Range("D38").Value = "'" + '6.75'
$ws.Range("E38").Value = '  -2.41%  '

$ws.Range("B39").Value = 'ImmutableX'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D39").Value = "'" + '1.50'
$ws.Range("E39").Value = '  -4.27%  '

$ws.Range("B40").Value = 'Monero'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D40").Value = "'" + '162.28'
$ws.Range("E40").Value = '  -1.00%  '

$ws.Range("D41").Value = "'" + '0.0777'
$ws.Range("E41").Value = '  -0.81%  '

$ws.Range("D42").Value = "'" + '0.798'
$ws.Range("E42").Value = '  -1.15%  '

$ws.Range("D43").Value = "'" + '0.999'
$ws.Range("E43").Value = '  -0.17%  '

$ws.Range("E44").Value = '  -1.05%  '

$ws.Range("D45").Value = "'" + '4.29'
$ws.Range("E45").Value = '  -2.93%  '

$ws.Range("D46").Value = "'" + '1.62'
$ws.Range("E46").Value = '  -2.13%  '

$ws.Range("B47").Value = 'ONDO'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D47").Value = "'" + '1.13'
$ws.Range("E47").Value = '  -3.23%  '

$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = "'" + '23.28'
$ws.Range("E48").Value = '  -7.47%  '

$ws.Range("D49").Value = "'" + '6.70'
$ws.Range("E49").Value = '  -1.24%  '

$ws.Range("D50").Value = "'" + '0.904'
$ws.Range("E50").Value = '  +0.78%  '

$ws.Range("D51").Value = "'" + '2.326.92'
$ws.Range("E51").Value = '  -5.42%  '
